# SE-1843: Aligned sample transaction/quote dates based on Future maturity dates
# Shift all sample date-strings from January 2021 to September 2021
# (keeping day-of-month and time-of-day the same), on both worksheets.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    $rows = $used.Rows.Count
    $cols = $used.Columns.Count
    $rowOffset = $used.Row
    $colOffset = $used.Column

    for ($r = 0; $r -lt $rows; $r++) {
        for ($c = 0; $c -lt $cols; $c++) {
            $cell = $ws.Cells.Item($r + $rowOffset, $c + $colOffset)
            $val = $cell.Value2
            if ($val -is [string] -and $val -like "2021-01-*") {
                $cell.Value2 = $val.Replace("2021-01-", "2021-09-")
            }
        }
    }
}

# Clear the stale selections left over from editing, matching the
# refreshed view state captured when the workbook was re-saved.
$ws1 = $wb.Worksheets.Item("transactions")
$ws2 = $wb.Worksheets.Item("prices")

$ws2.Range("A22").Select()
$ws1.Activate()
$ws1.Range("A1").Select()
